# Auto-generated edit script: applies scheduled-runner price/profit updates
# across the Leve-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 45318.082
$ws.Range("I62").Value = 74077.07000000001
$ws.Range("K62").Value = 74077.07000000001
$ws.Range("M62").Value = -73453.07000000001
$ws.Range("H64").Value = 4966.857
$ws.Range("I64").Value = 4121.6665
$ws.Range("J64").Value = 5600.75
$ws.Range("K64").Value = 4121.6665
$ws.Range("L64").Value = 5600.75
$ws.Range("M64").Value = -3873.6665
$ws.Range("N64").Value = -6096.75
$ws.Range("H65").Value = 45318.082
$ws.Range("I65").Value = 74077.07000000001
$ws.Range("K65").Value = 370385.35
$ws.Range("M65").Value = -367265.35
$ws.Range("H67").Value = 4966.857
$ws.Range("I67").Value = 4121.6665
$ws.Range("J67").Value = 5600.75
$ws.Range("K67").Value = 4121.6665
$ws.Range("L67").Value = 5600.75
$ws.Range("M67").Value = -3263.6665
$ws.Range("N67").Value = -7316.75
$ws.Range("H129").Value = 1323747.2
$ws.Range("I129").Value = 319.85715
$ws.Range("J129").Value = 1764889.6
$ws.Range("K129").Value = 959.5714499999999
$ws.Range("L129").Value = 5294668.800000001
$ws.Range("M129").Value = 4040.42855
$ws.Range("N129").Value = -5304668.800000001
$ws.Range("H137").Value = 947.625
$ws.Range("I137").Value = 940.1429000000001
$ws.Range("K137").Value = 2820.4287
$ws.Range("M137").Value = -270.4287000000004
$ws.Range("H138").Value = 3278.75
$ws.Range("I138").Value = 782.8946999999999
$ws.Range("J138").Value = 3966.0144
$ws.Range("K138").Value = 2348.6841
$ws.Range("L138").Value = 11898.0432
$ws.Range("M138").Value = 2791.3159
$ws.Range("N138").Value = -22178.0432
$ws.Range("H139").Value = 47306.25
$ws.Range("J139").Value = 47306.25
$ws.Range("L139").Value = 47306.25
$ws.Range("N139").Value = -57586.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H122").Value = 1205.6
$ws.Range("I122").Value = 1049.5
$ws.Range("K122").Value = 3148.5
$ws.Range("M122").Value = -698.5
$ws.Range("H139").Value = 42470
$ws.Range("J139").Value = 42470
$ws.Range("L139").Value = 42470
$ws.Range("N139").Value = -52750

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 83384300
$ws.Range("J20").Value = 166667680
$ws.Range("L20").Value = 166667680
$ws.Range("N20").Value = -166668174
$ws.Range("H137").Value = 40627.832
$ws.Range("J137").Value = 40627.832
$ws.Range("L137").Value = 40627.832
$ws.Range("N137").Value = -50827.832

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3146.6667
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 3366.6667
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 3366.6667
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -5612.6667
$ws.Range("H89").Value = 3146.6667
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 3366.6667
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 16833.3335
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -28065.3335
$ws.Range("H94").Value = 3055.875
$ws.Range("I94").Value = 2500
$ws.Range("J94").Value = 3135.2856
$ws.Range("K94").Value = 2500
$ws.Range("L94").Value = 3135.2856
$ws.Range("M94").Value = -2049
$ws.Range("N94").Value = -4037.2856
$ws.Range("H134").Value = 2779.5417
$ws.Range("I134").Value = 1845.45
$ws.Range("K134").Value = 5536.35
$ws.Range("M134").Value = -3001.35

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 167.5
$ws.Range("I33").Value = 140
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 840
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -557
$ws.Range("N33").Value = -2066
$ws.Range("H64").Value = 250577.75
$ws.Range("I64").Value = 770.6667
$ws.Range("J64").Value = 999999
$ws.Range("K64").Value = 2312.0001
$ws.Range("L64").Value = 2999997
$ws.Range("M64").Value = -2042.0001
$ws.Range("N64").Value = -3000537
$ws.Range("H67").Value = 250577.75
$ws.Range("I67").Value = 770.6667
$ws.Range("J67").Value = 999999
$ws.Range("K67").Value = 2312.0001
$ws.Range("L67").Value = 2999997
$ws.Range("M67").Value = -1376.0001
$ws.Range("N67").Value = -3001869
$ws.Range("H68").Value = 960.4
$ws.Range("I68").Value = 1025.5
$ws.Range("K68").Value = 3076.5
$ws.Range("M68").Value = -2265.5
$ws.Range("H71").Value = 960.4
$ws.Range("I71").Value = 1025.5
$ws.Range("K71").Value = 9229.5
$ws.Range("M71").Value = -5173.5
$ws.Range("H76").Value = 3013
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3013
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H82").Value = 9750
$ws.Range("J82").Value = 14500
$ws.Range("L82").Value = 43500
$ws.Range("N82").Value = -44312
$ws.Range("H85").Value = 9750
$ws.Range("J85").Value = 14500
$ws.Range("L85").Value = 43500
$ws.Range("N85").Value = -46308
$ws.Range("H107").Value = 3890388.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 3890388.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 11671165.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -11675005.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5316.6665
$ws.Range("H73").Value = 5316.6665
$ws.Range("H138").Value = 36733.332
$ws.Range("J138").Value = 37600
$ws.Range("L138").Value = 37600
$ws.Range("N138").Value = -47880

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 29000
$ws.Range("J62").Value = 29000
$ws.Range("L62").Value = 29000
$ws.Range("N62").Value = -30248
$ws.Range("H65").Value = 29000
$ws.Range("J65").Value = 29000
$ws.Range("L65").Value = 87000
$ws.Range("N65").Value = -93240
$ws.Range("H82").Value = 3429.5
$ws.Range("I82").Value = 3161.875
$ws.Range("J82").Value = 4500
$ws.Range("K82").Value = 3161.875
$ws.Range("L82").Value = 4500
$ws.Range("M82").Value = -2800.875
$ws.Range("N82").Value = -5222
$ws.Range("H85").Value = 3429.5
$ws.Range("I85").Value = 3161.875
$ws.Range("J85").Value = 4500
$ws.Range("K85").Value = 3161.875
$ws.Range("L85").Value = 4500
$ws.Range("M85").Value = -1913.875
$ws.Range("N85").Value = -6996
$ws.Range("H93").Value = 1682.7646
$ws.Range("I93").Value = 1016.9167
$ws.Range("J93").Value = 3280.8
$ws.Range("K93").Value = 1016.9167
$ws.Range("L93").Value = 3280.8
$ws.Range("M93").Value = 231.0833
$ws.Range("N93").Value = -5776.8
$ws.Range("H100").Value = 1145.75
$ws.Range("I100").Value = 1194.238
$ws.Range("J100").Value = 1000.2857
$ws.Range("K100").Value = 1194.238
$ws.Range("L100").Value = 1000.2857
$ws.Range("M100").Value = -653.2380000000001
$ws.Range("N100").Value = -2082.2857
$ws.Range("H134").Value = 38855.4
$ws.Range("J134").Value = 38855.4
$ws.Range("L134").Value = 38855.4
$ws.Range("N134").Value = -48995.4
$ws.Range("H138").Value = 32276.334
$ws.Range("J138").Value = 32276.334
$ws.Range("L138").Value = 32276.334
$ws.Range("N138").Value = -42556.334

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 26055.572
$ws.Range("J93").Value = 26055.572
$ws.Range("L93").Value = 26055.572
$ws.Range("N93").Value = -31047.572
